$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the week's hours (C10) and activity description (D10)
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = "investigacion html, css para aplicación de pagina en horizontal"

# Move the active selection to C12 (matches the author's last cursor position)
$ws.Range("C12").Select()
